$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.637.10'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '3.754.10'
$ws.Range("E3").Value = '  +6.79%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.91%  '
$ws.Range("D7").Value = '3.749.31'
$ws.Range("E7").Value = '  +6.34%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.538'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.166'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.34'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.493'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000253'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = '4.383.10'
$ws.Range("E15").Value = '  +7.04%  '
$ws.Range("D16").Value = '3.752.17'
$ws.Range("E16").Value = '  +7.00%  '
$ws.Range("D17").Value = '69.731.79'
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '513.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.07%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +17.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  +5.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.11%  '
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.52%  '
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.337'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.16'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.97%  '
$ws.Range("E40").Value = '  +5.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("E42").Value = '  -4.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '421.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.43%  '
$ws.Range("D45").Value = '3.074.62'
$ws.Range("E45").Value = '  +2.70%  '
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0363'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("E51").Value = '  -0.05%  '
